$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Contest 28 (row 37, "RR vs SRH") results
$ws.Range("E37").Value = 20
$ws.Range("H37").Value = 60
$ws.Range("K37").Value = 100
$ws.Range("N37").Value = 80
$ws.Range("Q37").Value = 0
$ws.Range("T37").Value = 40

# Contest 29 (row 38, "PBKS vs DC") results
$ws.Range("E38").Value = 20
$ws.Range("H38").Value = 100
$ws.Range("K38").Value = 60
$ws.Range("N38").Value = 80
$ws.Range("Q38").Value = 40
$ws.Range("T38").Value = 0
